$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Unit" values in column D (rows 20-23)
$ws.Range("D20").Value = "% of total"
$ws.Range("D21").Value = "% of total, SA"
$ws.Range("D22").Value = "% of GDP"
$ws.Range("D23").Value = "% of GDP, SA"

# New cell L20 referencing "% of GDP" string, with custom font styling
$ws.Range("L20").Value = "% of GDP"
$ws.Range("L20").Font.Name = "Courier New"
$ws.Range("L20").Font.Size = 11
$ws.Range("L20").Font.Color = 15323853

# Update selection / view
$ws.Range("E22").Select()
